$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '95.032.39'
$ws.Range("D3").Value = '3.565.66'
$ws.Range("E3").Value = '  -1.88%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.18'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.83%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '655.35'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.37%  '
$ws.Range("E7").Value = '  -1.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.400'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.45%  '
$ws.Range("E9").Value = '  +0.09%  '
$ws.Range("E10").Value = '  -0.04%  '
$ws.Range("D11").Value = '3.566.37'
$ws.Range("E11").Value = '  -1.73%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '42.44'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.44%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.45'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.43%  '
$ws.Range("D15").Value = '4.227.15'
$ws.Range("E15").Value = '  -2.17%  '
$ws.Range("D16").Value = '94.963.38'
$ws.Range("E16").Value = '  -1.37%  '
$ws.Range("E17").Value = '  -0.52%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.04'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.42%  '
$ws.Range("D19").Value = '3.543.87'
$ws.Range("E19").Value = '  -2.41%  '
$ws.Range("E20").Value = '  -4.94%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.83'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.78%  '
$ws.Range("E22").Value = '  +0.59%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '508.59'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.67%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.481'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.62%  '
$ws.Range("E25").Value = '  +3.61%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000197'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.60%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '91.69'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.95%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.74'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.85%  '
$ws.Range("D29").Value = '3.755.15'
$ws.Range("E29").Value = '  -1.76%  '
$ws.Range("E30").Value = '  -3.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.145'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.99%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.57'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.31%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.54%  '
$ws.Range("E35").Value = '  -3.12%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '31.93'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.38%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.73'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +16.85%  '
$ws.Range("B38").Value = 'Bittensor'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '604.09'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.59%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.59'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +9.13%  '
$ws.Range("E40").Value = '  -2.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.152'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.52%  '
$ws.Range("E42").Value = '  +0.17%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.906'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.85'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +6.42%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '35.20'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +28.18%  '
$ws.Range("E46").Value = '  +0.35%  '
$ws.Range("E47").Value = '  +3.52%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.40'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.70%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.49'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.54%  '
$ws.Range("E51").Value = '  -0.08%  '
